# Convert the single "login" test-data sheet into two sheets:
#   - ValidLogin   : the existing valid username/password pair
#   - InvalidLogin : a new sheet with an invalid username/password pair
# (commit message: "Invalid login test case")

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Rework the first sheet -------------------------------------------------
$ws1.Name = "ValidLogin"

# Wipe the old 3-column / 4-row table (Username/Password/Status + 3 data rows)
$ws1.UsedRange.Clear()

# Rebuild it as a simple 2-column UserName/Password table
$ws1.Range("A1").Value = "UserName"
$ws1.Range("B1").Value = "Password"
$ws1.Range("A2").Value = "admin"
$ws1.Range("B2").Value = "manager"

# --- Add the new sheet for the invalid-login test case ----------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "InvalidLogin"

$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "abcd"
$ws2.Range("B2").Value = "xasdas"

# --- Selections / active sheet ----------------------------------------------
$ws1.Range("A1:B2").Select()

$ws2.Activate()
$ws2.Range("B3").Select()
